$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New objective (B), gap (C), and solve time (D) values for rows 2-11
$data = @(
    @{Row=2;  B=-105.20922885539122; C=0.09054680418354322; D=125.264838653},
    @{Row=3;  B=-100.25578112906734; C=0.06162279836591759; D=134.873800722},
    @{Row=4;  B=-103.67140169613268; C=0.02259346017021363; D=210.849153706},
    @{Row=5;  B=-103.09145033800306; C=0.09806946180741519; D=70.472930477},
    @{Row=6;  B=-102.05315995415073; C=0.09107443008400061; D=87.859177768},
    @{Row=7;  B=-102.01057739411644; C=0.0957130934643146;  D=150.70766555},
    @{Row=8;  B=-97.44343274822958;  C=0.006862235569463699; D=12.982388247},
    @{Row=9;  B=-102.53029754612697; C=0.09900204050256638; D=122.121381636},
    @{Row=10; B=-102.29121992425158; C=0.06533247987843727; D=43.570730038},
    @{Row=11; B=-99.01915304554905;  C=0.08516894509857384; D=107.592718432}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}
